$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I (I0) and J (IF) -- match the bold/centered
# header style already used by H1 (and the rest of row 1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy column H values (IP) into the new J column (IF), and default
# column I (I0) to 1, for every data row 2..34.
for ($r = 2; $r -le 34; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 32 is a special case where I0/IF diverge from the default pattern.
$ws.Cells.Item(32, 9).Value = 5
$ws.Cells.Item(32, 10).Value = 8
